# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund holdings detail) positioned
#    right after "2021-Q4" and before "总计".
# 2. Populate it with the same column layout used by the other quarterly
#    sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# 3. Insert a new summary row at the top of "总计" for 2022-Q1 and bump
#    the index column of the existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create and place the new "2022-Q1" worksheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# --- header row -------------------------------------------------------
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Value = "基金代码"
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- row 2 : 513500 博时标普500ETF(QDII) --------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("A2").Font.Bold = $true
$newSheet.Range("A2").HorizontalAlignment = -4108
$newSheet.Range("A2").VerticalAlignment = -4160
$newSheet.Range("A2").Borders.LineStyle = 1

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "513500"
$newSheet.Range("B2").Style = "Normal"

$newSheet.Range("C2").NumberFormat = "@"
$newSheet.Range("C2").Value = "博时标普500ETF(QDII)"
$newSheet.Range("C2").Style = "Normal"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "70.03"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "90.45"
$newSheet.Range("E2").Style = "Normal"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "1.53"
$newSheet.Range("F2").Style = "Normal"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "1.0715"
$newSheet.Range("G2").Style = "Normal"

$newSheet.Range("H2").Value = 8

# --- row 3 : 003718 易方达标普500指数(QDII-LOF) 美元 ----------------------
$newSheet.Range("A3").Value = 1
$newSheet.Range("A3").Font.Bold = $true
$newSheet.Range("A3").HorizontalAlignment = -4108
$newSheet.Range("A3").VerticalAlignment = -4160
$newSheet.Range("A3").Borders.LineStyle = 1

$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "003718"
$newSheet.Range("B3").Style = "Normal"

$newSheet.Range("C3").NumberFormat = "@"
$newSheet.Range("C3").Value = "易方达标普500指数(QDII-LOF) 美元"
$newSheet.Range("C3").Style = "Normal"

$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "5.22"
$newSheet.Range("D3").Style = "Normal"

$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "91.11"
$newSheet.Range("E3").Style = "Normal"

$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "1.53"
$newSheet.Range("F3").Style = "Normal"

$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0799"
$newSheet.Range("G3").Style = "Normal"

$newSheet.Range("H3").Value = 8

# --- row 4 : 161125 易方达标普500指数(QDII-LOF) 人民币 --------------------
$newSheet.Range("A4").Value = 2
$newSheet.Range("A4").Font.Bold = $true
$newSheet.Range("A4").HorizontalAlignment = -4108
$newSheet.Range("A4").VerticalAlignment = -4160
$newSheet.Range("A4").Borders.LineStyle = 1

$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B4").Value = "161125"
$newSheet.Range("B4").Style = "Normal"

$newSheet.Range("C4").NumberFormat = "@"
$newSheet.Range("C4").Value = "易方达标普500指数(QDII-LOF) 人民币"
$newSheet.Range("C4").Style = "Normal"

$newSheet.Range("D4").NumberFormat = "@"
$newSheet.Range("D4").Value = "5.22"
$newSheet.Range("D4").Style = "Normal"

$newSheet.Range("E4").NumberFormat = "@"
$newSheet.Range("E4").Value = "91.11"
$newSheet.Range("E4").Style = "Normal"

$newSheet.Range("F4").NumberFormat = "@"
$newSheet.Range("F4").Value = "1.53"
$newSheet.Range("F4").Style = "Normal"

$newSheet.Range("G4").NumberFormat = "@"
$newSheet.Range("G4").Value = "0.0799"
$newSheet.Range("G4").Style = "Normal"

$newSheet.Range("H4").Value = 8

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new 2022-Q1 row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
# the freshly-inserted row inherits stray formatting from the header
# row above it; strip that back to Normal before (re)building it
$total.Range("A2:D2").Style = "Normal"

$total.Range("A2").Value = 0
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

$total.Range("B2").NumberFormat = "@"
$total.Range("B2").Value = "2022-Q1"
$total.Range("B2").Style = "Normal"

$total.Range("C2").Value = 3
$total.Range("D2").Value = 1.23

# bump the index column for the two pre-existing rows, now shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

Write-Host "edit complete"
